# The "Förändrad" (changed/modified) date in column C was bumped by one day
# (2023-10-05 -> 2023-10-06, serial 45204 -> 45205) for every data row
# (rows 2 through 163) on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45205
    }
}
